$d = $word.ActiveDocument

$pairs = @(
    @("663÷4=165, 3", "443÷8=55, 3"),
    @("896÷8=112, 0", "971÷8=121, 3"),
    @("716÷2=358, 0", "619÷9=68, 7"),
    @("333÷5=66, 3", "779÷8=97, 3"),
    @("665÷6=110, 5", "810÷6=135, 0"),
    @("344÷7=49, 1", "230÷8=28, 6"),
    @("439÷6=73, 1", "570÷9=63, 3"),
    @("284÷8=35, 4", "514÷4=128, 2"),
    @("161÷3=53, 2", "524÷5=104, 4"),
    @("641÷8=80, 1", "400÷4=100, 0"),
    @("518÷2=259, 0", "425÷4=106, 1"),
    @("400÷2=200, 0", "361÷9=40, 1"),
    @("582÷5=116, 2", "950÷8=118, 6"),
    @("704÷6=117, 2", "676÷6=112, 4"),
    @("722÷3=240, 2", "143÷7=20, 3"),
    @("687÷6=114, 3", "546÷6=91, 0"),
    @("557÷9=61, 8", "366÷4=91, 2"),
    @("495÷3=165, 0", "861÷7=123, 0"),
    @("507÷3=169, 0", "454÷8=56, 6"),
    @("163÷8=20, 3", "824÷3=274, 2"),
    @("297÷2=148, 1", "610÷5=122, 0"),
    @("409÷6=68, 1", "838÷7=119, 5"),
    @("626÷6=104, 2", "322÷8=40, 2"),
    @("261÷2=130, 1", "362÷3=120, 2"),
    @("750÷3=250, 0", "906÷2=453, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
